# Applies the edit described by the commit:
# "Tested with 1 camera (multi camera opmode unreliable)"
# - Adds a new worksheet "Sheet5" (a re-run of the id1 Straight-On test,
#   using only rows for tags 0-0 .. 0-9, with tag 0-1 swapped out for a
#   second 0-4 reading) and updates the selections on several sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("id1 - Straight On (Angle Cxled)")
$ws2 = $wb.Worksheets.Item("id2 - Pointing")
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Create the new "Sheet5" worksheet as the last sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet5"

# --- Header row ---
$newSheet.Range("A1").Value = "Location"
$newSheet.Range("D1").Value = "Predicted Angle"
$newSheet.Range("E1").Value = "Predicted X"
$newSheet.Range("F1").Value = "Predicted Y"
$newSheet.Range("G1").Value = "Actual Angle"
$newSheet.Range("H1").Value = "Actual X"
$newSheet.Range("I1").Value = "Actual Y"

# --- Data rows (rows 2-11) ---
$newSheet.Range("A2").Value = "0-0"
$newSheet.Range("D2").Value = -0.63
$newSheet.Range("E2").Value = 5.84
$newSheet.Range("F2").Value = 23.92
$newSheet.Range("G2").Value = 0
$newSheet.Range("H2").Value = 6
$newSheet.Range("I2").Value = 24

$newSheet.Range("A3").Value = "0-4"
$newSheet.Range("D3").Value = -0.64
$newSheet.Range("E3").Value = 5.94
$newSheet.Range("F3").Value = 16.079999999999998
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 0
$newSheet.Range("I3").Value = 16

$newSheet.Range("A4").Value = "0-2"
$newSheet.Range("D4").Value = -0.89
$newSheet.Range("E4").Value = 0.32
$newSheet.Range("F4").Value = 20.350000000000001
$newSheet.Range("G4").Value = 0
$newSheet.Range("H4").Value = 0
$newSheet.Range("I4").Value = 20

$newSheet.Range("A5").Value = "0-3"
$newSheet.Range("D5").Value = -1.1000000000000001
$newSheet.Range("E5").Value = 0.35
$newSheet.Range("F5").Value = 18.27
$newSheet.Range("G5").Value = 0
$newSheet.Range("H5").Value = 0
$newSheet.Range("I5").Value = 18

$newSheet.Range("A6").Value = "0-4"
$newSheet.Range("D6").Value = -1.18
$newSheet.Range("E6").Value = 0.34
$newSheet.Range("F6").Value = 16.510000000000002
$newSheet.Range("G6").Value = 0
$newSheet.Range("H6").Value = 0
$newSheet.Range("I6").Value = 16

$newSheet.Range("A7").Value = "0-5"
$newSheet.Range("D7").Value = -0.85
$newSheet.Range("E7").Value = 0.22
$newSheet.Range("F7").Value = 14.58
$newSheet.Range("G7").Value = 0
$newSheet.Range("H7").Value = 0
$newSheet.Range("I7").Value = 14

$newSheet.Range("A8").Value = "0-6"
$newSheet.Range("D8").Value = 0.32
$newSheet.Range("E8").Value = -0.07
$newSheet.Range("F8").Value = 12.57
$newSheet.Range("G8").Value = 0
$newSheet.Range("H8").Value = 0
$newSheet.Range("I8").Value = 12

$newSheet.Range("A9").Value = "0-7"
$newSheet.Range("D9").Value = 0.08
$newSheet.Range("E9").Value = -0.01
$newSheet.Range("F9").Value = 10.66
$newSheet.Range("G9").Value = 0
$newSheet.Range("H9").Value = 0
$newSheet.Range("I9").Value = 10

$newSheet.Range("A10").Value = "0-8"
$newSheet.Range("D10").Value = 0.23
$newSheet.Range("E10").Value = -0.03
$newSheet.Range("F10").Value = 8.66
$newSheet.Range("G10").Value = 0
$newSheet.Range("H10").Value = 0
$newSheet.Range("I10").Value = 8

$newSheet.Range("A11").Value = "0-9"
$newSheet.Range("D11").Value = -0.17
$newSheet.Range("E11").Value = 0.02
$newSheet.Range("F11").Value = 6.53
$newSheet.Range("G11").Value = 0
$newSheet.Range("H11").Value = 0
$newSheet.Range("I11").Value = 6

# --- Apply the same ("Menlo" 11pt) cell style used by Sheet1's data cells ---
# (copy cell-formatting only, from a source cell on Sheet1 already using
# that style, onto every target cell that needs it)
$styleSource = $ws1.Range("D3")
$styleSource.Copy()
$styledCells = @("E2","F2","D3","E3","D4","E4","F4","D5","E5","F5","D6","E6","F6","D7","E7","F7","D8","E8","F8","D9","E9","F9","D10","E10","F10","D11","E11","F11")
foreach ($cellRef in $styledCells) {
    $newSheet.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Selection / active-cell bookkeeping ---

# Sheet1: selection grew from the header row to the whole data block
$ws1.Activate()
$ws1.Range("A1:I11").Select()

# Sheet2: active cell moved to A2
$ws2.Activate()
$ws2.Range("A2").Select()

# Sheet4: no longer the active tab; active cell moved to B1
$ws4.Activate()
$ws4.Range("B1").Select()

# Sheet5: newly added sheet becomes the active tab, active cell A4
$newSheet.Activate()
$newSheet.Range("A4").Select()
